# Apply the update: add two new medicine rows (BONDIGA 120ML SYRUP and
# DOLIPRANE 1 GM 15 TABS.) at the top of the data table, renumber the
# existing rows, and refresh the totals row with the new sum.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# 1) Insert two blank rows right above the existing first data row (row 4).
#    This pushes all the existing data rows, the totals row and the footer
#    row down by two rows, and Excel keeps the existing merged cells and
#    row heights intact for everything below.
$ws.Range("A4:N5").EntireRow.Insert()

# 2) Copy the formatting (styles/borders/fonts/fills) of a normal data row
#    (now at row 6, the old row 4) onto the two freshly inserted rows so
#    they look like the rest of the table.
$ws.Range("A6:N6").Copy()
$ws.Range("A4:N4").PasteSpecial($xlPasteFormats)
$ws.Range("A6:N6").Copy()
$ws.Range("A5:N5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# 3) Restore the exact row heights used in the final layout.
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 24.75

# 4) Re-create the merged cell groups for the two new rows, matching the
#    pattern used by every other data row (name / ratio / price columns).
$ws.Range("B4:G4").Merge()
$ws.Range("H4:K4").Merge()
$ws.Range("L4:M4").Merge()
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()

# 5) Populate the two new data rows.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "BONDIGA 120ML  SYRUP"
$ws.Cells.Item(4, 8).Value = "0:0"
$ws.Cells.Item(4, 12).Value = 69
$ws.Cells.Item(4, 14).Value = "1:0"

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "DOLIPRANE 1 GM 15 TABS."
$ws.Cells.Item(5, 8).Value = "7:1"
$ws.Cells.Item(5, 12).Value = 16
$ws.Cells.Item(5, 14).Value = "0:0"

# 6) Renumber the "م" (index) column for the rows that shifted down.
$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(11, 1).Value = 8

# 7) Update the grand total (row 12, shifted down from row 10) with the new
#    sum of the price column.
$ws.Cells.Item(12, 11).Value = 414.5

Write-Host "Edit applied successfully"
